$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")
$ws.Range("L20").Formula = "=15096.928674"
